$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows M:T for rows 2-4 with refreshed TPM values

# Row 2
$ws.Range("M2").Value = 19.524618
$ws.Range("N2").Value = 58.573854
$ws.Range("O2").Value = 0.4154885426712971
$ws.Range("P2").Value = 0.4539723485554654
$ws.Range("Q2").Value = 7.583387664024001
$ws.Range("R2").Value = 68.250488976216
$ws.Range("S2").Value = 0.4154885426712971
$ws.Range("T2").Value = 0.4539723485554654

# Row 3
$ws.Range("M3").Value = 15.24435933333334
$ws.Range("N3").Value = 45.73307800000001
$ws.Range("O3").Value = 0.324403614112412
$ws.Range("P3").Value = 0.3544508583357054
$ws.Range("Q3").Value = 5.920929490879113
$ws.Range("R3").Value = 53.28836541791201
$ws.Range("S3").Value = 0.324403614112412
$ws.Range("T3").Value = 0.3544508583357054

# Row 4
$ws.Range("M4").Value = 11.9507005
$ws.Range("N4").Value = 23.901401
$ws.Range("O4").Value = 0.2543137660693869
$ws.Range("P4").Value = 0.1852460510065796
$ws.Range("Q4").Value = 4.641668008467334
$ws.Range("R4").Value = 27.850008050804
$ws.Range("S4").Value = 0.2543137660693869
$ws.Range("T4").Value = 0.1852460510065796

# New row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3884013333333334
$ws.Range("H5").Value = 1.165204
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.272275
$ws.Range("N5").Value = 0.816825
$ws.Range("O5").Value = 0.005794077146903843
$ws.Range("P5").Value = 0.006330742102249548
$ws.Range("Q5").Value = 0.1057519730333333
$ws.Range("R5").Value = 0.9517677573000002
$ws.Range("S5").Value = 0.005794077146903843
$ws.Range("T5").Value = 0.006330742102249548
